$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the top. This shifts the existing header
# row (old row 1 -> row 3), data row (old row 2 -> row 4), the
# "Use text format..." note (old row 3 -> row 5) and the two trailing
# spacer rows (old row 31 -> 33, old row 58 -> 60) down by two rows.
$ws.Rows("1:2").Insert()

# The engine does not re-anchor existing Hyperlink objects when rows are
# inserted above them, so drop the stale one (still pointing at E2) and
# recreate it at its new location E4, preserving the original target URL.
$ws.Hyperlinks.Delete()
$hlCell = $ws.Range("E4")
$ws.Hyperlinks.Add($hlCell, "https://more.arrs.run/runner/15799 ") | Out-Null
$hlCell.NumberFormat = "@"

# New title cell explaining the purpose of this sheet.
$ws.Range("A1").Value2 = "Add known performances here (e.g. for road events) that predate powerof10 and are missing from the official C&C T&F records"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Size = 14
$ws.Range("A1").Font.Name = "Arial"

# Update the notes/comments that previously questioned whether the
# performance was a club PB - now confirmed.
$ws.Range("G4").Value2 = "Harlow"
$ws.Range("K4").Value2 = "This matches 2009 club records so must be right performance"

# Restore the active selection to the note cell under the table.
$ws.Range("A5").Select()
